# Applies the "time-log" update described in the commit:
#   - Logs two new entries (rows 15 and 16) with their Notes text
#   - Adds Notes text to rows 12 and 13 (B/C/D already had data there)
#   - Extends the shared "duration" formula down through row 16
#   - Extends the Total Time sum to include the new rows
#   - Updates the active-cell selection to E16
#   - Nudges the saved window position (xWindow) back to 0

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New shared strings must be entered in the same order the author typed them so
#     the shared-string table indices line up with the target workbook ---

# Row 15 note first (becomes shared string index 16)
$ws.Range("E15").Value = "Cleaning and preparing data files for sending to client."

# Notes for the already-populated rows 12 & 13 (indices 17 & 18)
$ws.Range("E12").Value = "Revising the input variables used in the principal components scores and updating groupings and clusters"
$ws.Range("E13").Value = "Editing visualizations and adding new figures to explore the distribution of the MSAs across different combinations of the principal component scores and cluster IDs"

# Row 16 note last (index 19)
$ws.Range("E16").Value = "Revising to 35 clusters. Also modifying the quantiles to try to get less strange cluster groupings."

# --- Row 15: new Start/End entries ---
$ws.Range("B15").Value = 0.43055555555555558
$ws.Range("C15").Value = 0.4604166666666667

# --- Row 16: new Start entry (no End) ---
$ws.Range("B16").Value = 0.3888888888888889

# --- Duration formula (C-B) extended down through the two new rows ---
$ws.Range("D15:D16").Formula = "=C15-B15"

# --- Extend the Total Time sum to include the new rows ---
$ws.Range("D20").Formula = "=SUM(D2:D16)"

# --- Update the saved selection / active cell ---
[void]$ws.Range("E16").Select()

# --- Restore window x-position (was shifted during editing, best effort) ---
$wb.Windows.Item(1).Left = 0
